$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for team record, matching the style of the existing
# header row (bold, centered, thin border - same as AC1 "Unnamed: 28").
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Every player row (2-48) gets the team's season record repeated across
# the new Wins / Losses / Ties columns.
for ($r = 2; $r -le 48; $r++) {
    $ws.Cells.Item($r, 30).Value = 78  # AD - Wins
    $ws.Cells.Item($r, 31).Value = 84  # AE - Losses
    $ws.Cells.Item($r, 32).Value = 0   # AF - Ties
}
